$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current header row (row 1). This pushes
# every existing row down by one (old row 1 -> row 2, old row 2 -> row 3, ...,
# old row 22 -> row 23).
$ws.Rows("1:1").Insert()

# The newly inserted row 1 is blank/unformatted and the old header row (now
# row 2) still carries the bold/bordered header style. Copy that formatting
# up to the new row 1 so the numeric index row gets the header look, then
# reset row 2 back to the default "Normal" style (it should just hold plain
# text now).
$ws.Range("A2:M2").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A2:M2").Style = "Normal"

# Fill the new row 1 with a simple numeric column index (0-based), replacing
# the old textual column headers.
for ($i = 0; $i -le 12; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# The old L1/M1 header labels ("thread_size" / "material_surface") shifted
# down into L2/M2 along with the rest of row 1 -- but those two no longer
# belong there, so clear them out.
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
